$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 106
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()
# Row 137
$ws.Range("H137").Value = 1984.25
$ws.Range("I137").Value = 1459.4
$ws.Range("K137").Value = 4378.200000000001
$ws.Range("M137").Value = -1828.200000000001
# Row 138
$ws.Range("H138").Value = 1672.8857
$ws.Range("I138").Value = 1253.3513
$ws.Range("J138").Value = 2143.2727
$ws.Range("K138").Value = 3760.0539
$ws.Range("L138").Value = 6429.8181
$ws.Range("M138").Value = 1379.9461
$ws.Range("N138").Value = -16709.8181

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 935.75
$ws.Range("I2").Value = 948.1579
$ws.Range("J2").Value = 700
$ws.Range("K2").Value = 948.1579
$ws.Range("L2").Value = 700
$ws.Range("M2").Value = -835.1579
$ws.Range("N2").Value = -926
# Row 61
$ws.Range("H61").Value = 1404.9697
$ws.Range("I61").Value = 1385.9375
$ws.Range("K61").Value = 1385.9375
$ws.Range("M61").Value = -1173.9375
# Row 88
$ws.Range("H88").Value = 2760
$ws.Range("I88").Value = 1300
$ws.Range("J88").Value = 3125
$ws.Range("K88").Value = 1300
$ws.Range("L88").Value = 3125
$ws.Range("M88").Value = -894
$ws.Range("N88").Value = -3937
# Row 91
$ws.Range("H91").Value = 2760
$ws.Range("I91").Value = 1300
$ws.Range("J91").Value = 3125
$ws.Range("K91").Value = 1300
$ws.Range("L91").Value = 3125
$ws.Range("M91").Value = 104
$ws.Range("N91").Value = -5933
# Row 116
$ws.Range("H116").Value = 935.75
$ws.Range("I116").Value = 948.1579
$ws.Range("J116").Value = 700
$ws.Range("K116").Value = 948.1579
$ws.Range("L116").Value = 700
$ws.Range("M116").Value = 1345.8421
$ws.Range("N116").Value = -5288
# Row 132
$ws.Range("H132").Value = 1914.7778
$ws.Range("I132").Value = 1048.25
$ws.Range("J132").Value = 3175.182
$ws.Range("K132").Value = 3144.75
$ws.Range("L132").Value = 9525.545999999998
$ws.Range("M132").Value = -614.75
$ws.Range("N132").Value = -14585.546
# Row 136
$ws.Range("H136").Value = 1404.9697
$ws.Range("I136").Value = 1385.9375
$ws.Range("K136").Value = 4157.8125
$ws.Range("M136").Value = -1607.8125

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 935.75
$ws.Range("I3").Value = 948.1579
$ws.Range("J3").Value = 700
$ws.Range("K3").Value = 948.1579
$ws.Range("L3").Value = 700
$ws.Range("M3").Value = -834.1579
$ws.Range("N3").Value = -928
# Row 105
$ws.Range("H105").Value = 3163.3157
$ws.Range("I105").Value = 2859
$ws.Range("J105").Value = 5750
$ws.Range("K105").Value = 2859
$ws.Range("L105").Value = 5750
$ws.Range("M105").Value = -1112
$ws.Range("N105").Value = -9244
# Row 107
$ws.Range("H107").Value = 1752.9642
$ws.Range("I107").Value = 1441.1666
$ws.Range("J107").Value = 3623.75
$ws.Range("K107").Value = 1441.1666
$ws.Range("L107").Value = 3623.75
$ws.Range("M107").Value = 478.8334
$ws.Range("N107").Value = -7463.75
# Row 134
$ws.Range("H134").Value = 2014.9678
$ws.Range("I134").Value = 1084.7826
$ws.Range("J134").Value = 4689.25
$ws.Range("K134").Value = 3254.3478
$ws.Range("L134").Value = 14067.75
$ws.Range("M134").Value = -719.3478
$ws.Range("N134").Value = -19137.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2325.2559
$ws.Range("I31").Value = 2036.45
$ws.Range("J31").Value = 2576.3914
$ws.Range("K31").Value = 2036.45
$ws.Range("L31").Value = 2576.3914
$ws.Range("M31").Value = -1741.45
$ws.Range("N31").Value = -3166.3914
# Row 34
$ws.Range("H34").Value = 2325.2559
$ws.Range("I34").Value = 2036.45
$ws.Range("J34").Value = 2576.3914
$ws.Range("K34").Value = 2036.45
$ws.Range("L34").Value = 2576.3914
$ws.Range("M34").Value = -1834.45
$ws.Range("N34").Value = -2980.3914
# Row 132
$ws.Range("H132").Value = 2843.2144
$ws.Range("I132").Value = 1373.1428
$ws.Range("J132").Value = 4313.2856
$ws.Range("K132").Value = 4119.428400000001
$ws.Range("L132").Value = 12939.8568
$ws.Range("M132").Value = -1589.428400000001
$ws.Range("N132").Value = -17999.8568

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1762.5454
$ws.Range("J5").Value = 6495
$ws.Range("L5").Value = 19485
$ws.Range("N5").Value = -19709
# Row 122
$ws.Range("H122").Value = 674.25
$ws.Range("I122").Value = 340
$ws.Range("J122").Value = 934.2222
$ws.Range("K122").Value = 3060
$ws.Range("L122").Value = 8407.9998
$ws.Range("M122").Value = -610
$ws.Range("N122").Value = -13307.9998
# Row 129
$ws.Range("H129").Value = 1774.1428
$ws.Range("I129").Value = 764.2857
$ws.Range("J129").Value = 2784
$ws.Range("K129").Value = 2292.8571
$ws.Range("L129").Value = 8352
$ws.Range("M129").Value = 2707.1429
$ws.Range("N129").Value = -18352
# Row 131
$ws.Range("H131").Value = 4628.7095
$ws.Range("J131").Value = 4920.6895
$ws.Range("L131").Value = 14762.0685
$ws.Range("N131").Value = -24842.0685
# Row 135
$ws.Range("H135").Value = 1762.5454
$ws.Range("J135").Value = 6495
$ws.Range("L135").Value = 58455
$ws.Range("N135").Value = -63525

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 1884
$ws.Range("I122").Value = 1663.4
$ws.Range("J122").Value = 2686.182
$ws.Range("K122").Value = 4990.200000000001
$ws.Range("L122").Value = 8058.545999999999
$ws.Range("M122").Value = -2540.200000000001
$ws.Range("N122").Value = -12958.546

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 87984.164
$ws.Range("I40").Value = 251000
$ws.Range("J40").Value = 6476.25
$ws.Range("K40").Value = 251000
$ws.Range("L40").Value = 6476.25
$ws.Range("M40").Value = -250864
$ws.Range("N40").Value = -6748.25
# Row 100
$ws.Range("H100").Value = 1770
$ws.Range("I100").Value = 1026.6666
$ws.Range("J100").Value = 4000
$ws.Range("K100").Value = 1026.6666
$ws.Range("L100").Value = 4000
$ws.Range("M100").Value = -485.6666
$ws.Range("N100").Value = -5082
# Row 132
$ws.Range("H132").Value = 6665.9434
$ws.Range("I132").Value = 7448.086
$ws.Range("J132").Value = 5145.1113
$ws.Range("K132").Value = 22344.258
$ws.Range("L132").Value = 15435.3339
$ws.Range("M132").Value = -19814.258
$ws.Range("N132").Value = -20495.3339

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 1316.7273
$ws.Range("I96").Value = 1272.5
$ws.Range("J96").Value = 1434.6666
$ws.Range("K96").Value = 1272.5
$ws.Range("L96").Value = 1434.6666
$ws.Range("M96").Value = 100.5
$ws.Range("N96").Value = -4180.6666
